$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet only needs to stay as a blank, documentation-only template now
# (runners / chefs will sign in their own orders later), so wipe out the
# sample order data that was filled in for testing. Row 2's Order # (A2)
# is left alone; everything else in the sample block (B2:H2 and A3:H5) is
# cleared back out to blank cells, matching the rest of the still-empty
# rows further down the sheet.
$ws.Range("B2:H2").ClearContents()
$ws.Range("A3:H5").ClearContents()
